$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "JUnit / Mockito / Cucumber / Gherkin",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "JUnit / Mockito / Cucumber / Gherkin / Wiremock",
    2
)

$d.Content.Find.Execute(
    "Git / Maven / Jenkins / JIRA",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Git / Maven / Jenkins / JIRA / SonarQube",
    2
)
